$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "assignment" value in F10 from "by_population" to "by_pop"
$ws.Range("F10").Value = "by_pop"

# Match the resulting selection/active cell from the edit
$ws.Range("F10").Select()
